$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.751.13"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.538.10"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Formula = "'311.16"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Formula = "'100.67"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").Formula = "'35.70"
$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Formula = "'7.32"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "2.931.00"

$ws.Range("D15").Value = "2.603.09"
$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("E16").Value = "  -2.90%  "

$ws.Range("D17").Formula = "'0.815"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "42.752.18"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Formula = "'70.14"

$ws.Range("D23").Formula = "'243.31"
$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("D27").Formula = "'25.46"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("D29").Formula = "'10.18"
$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("D30").Formula = "'38.66"
$ws.Range("E30").Value = "  -4.20%  "

$ws.Range("D31").Formula = "'159.21"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("E33").Value = "  +6.18%  "

$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("D35").Formula = "'0.0792"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Formula = "'18.14"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Formula = "'3.15"
$ws.Range("E37").Value = "  -3.76%  "

$ws.Range("E38").Value = "  -5.09%  "

$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").Formula = "'4.13"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("D42").Formula = "'21.78"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").Value = "1.995.57"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").Formula = "'9.15"
$ws.Range("E47").Value = "  +1.19%  "

$ws.Range("D48").Value = "2.782.30"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").Formula = "'0.192"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Formula = "'80.01"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("D51").Formula = "'72.31"
$ws.Range("E51").Value = "  -1.23%  "
